$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1931
$ws.Range("F3").Value = 809
$ws.Range("F4").Value = 13531
$ws.Range("F5").Value = 13353
$ws.Range("F6").Value = 1033
$ws.Range("F9").Value = 579
$ws.Range("F12").Value = 18
$ws.Range("F13").Value = 713
$ws.Range("F14").Value = 2117
$ws.Range("F15").Value = 39
$ws.Range("F16").Value = 74
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 60
$ws.Range("F18").Value = 94
$ws.Range("F21").Value = 322
$ws.Range("F23").Value = 469
$ws.Range("F24").Value = 793
$ws.Range("F25").Value = 51

$ws = $wb.Worksheets.Item(2)
$ws.Range("G2").Value = "不可售"
$ws.Range("F6").Value = 67
$ws.Range("F8").Value = 823

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 207
$ws.Range("F3").Value = 78

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 207
$ws.Range("F3").Value = 1931
$ws.Range("G4").Value = "不可售"
$ws.Range("F5").Value = 809
$ws.Range("F6").Value = 13531
$ws.Range("F7").Value = 13353
$ws.Range("F8").Value = 1033
$ws.Range("F11").Value = 579
$ws.Range("F14").Value = 18
$ws.Range("F15").Value = 713
$ws.Range("F18").Value = 2117
$ws.Range("F19").Value = 39
$ws.Range("F20").Value = 74
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 60
$ws.Range("F22").Value = 94
$ws.Range("F25").Value = 67
$ws.Range("F26").Value = 78
$ws.Range("F28").Value = 322
$ws.Range("F30").Value = 469
$ws.Range("F31").Value = 793
$ws.Range("F33").Value = 823
$ws.Range("F36").Value = 51
